$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (Priyanka Muddana) ---
$ws.Range("C2").Value = "Priyanka Muddana"

# D2 holds a date-like string ("02/06/2014") that must stay plain text,
# not get auto-converted to a date serial by Excel's smart-entry parsing.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "02/06/2014"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "Internship"
$ws.Range("F2").Value = "HR Manager"

$ws.Range("H2").Value = 120000.0
$ws.Range("I2").Value = 10000.0
$ws.Range("J2").Value = 10345.0
$ws.Range("K2").Value = 30.0
$ws.Range("L2").Value = 30.0
$ws.Range("M2").Value = 4000.0
$ws.Range("N2").Value = 1000.0
$ws.Range("O2").Value = 4345.0
$ws.Range("P2").Value = 1000.0
$ws.Range("Q2").Value = 10345.0
$ws.Range("R2").Value = 480.0
$ws.Range("S2").Value = 181.04
$ws.Range("T2").Value = 150.0
$ws.Range("U2").Value = 155.0
$ws.Range("V2").Value = 100.0
$ws.Range("W2").Value = 1066.04
$ws.Range("X2").Value = 9278.96

# --- Add new row 3 (Vidya Sagar pogiri) ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "11-2014"
$ws.Range("C3").Value = "Vidya Sagar pogiri"

# D3 holds a date-like string too; keep it plain text the same way.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "02/06/2014"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "Regular"
$ws.Range("F3").Value = "Junior Developer"
$ws.Range("G3").Value = "Development"

$ws.Range("H3").Value = 130000.0
$ws.Range("I3").Value = 10833.333333333334
$ws.Range("J3").Value = 10313.3
$ws.Range("K3").Value = 30.0
$ws.Range("L3").Value = 30.0
$ws.Range("M3").Value = 4333.33
$ws.Range("N3").Value = 1083.33
$ws.Range("O3").Value = 4896.67
$ws.Range("P3").Value = 0.0
$ws.Range("Q3").Value = 10313.3
$ws.Range("R3").Value = 520.0
$ws.Range("S3").Value = 0.0
$ws.Range("T3").Value = 0.0
$ws.Range("U3").Value = 0.0
$ws.Range("V3").Value = 0.0
$ws.Range("W3").Value = 520.0
$ws.Range("X3").Value = 9793.33

# --- Column width adjustments (auto-fit side effect of the longer content
#     now in these columns - widths match the workbook's post-edit layout) ---
$ws.Columns.Item(3).ColumnWidth = 13.358301926163726
$ws.Columns.Item(4).ColumnWidth = 13.358301926163726
$ws.Columns.Item(6).ColumnWidth = 14.455601926163727
$ws.Columns.Item(16).ColumnWidth = 7.855601926163724
$ws.Columns.Item(18).ColumnWidth = 6.755601926163725
$ws.Columns.Item(20).ColumnWidth = 6.755601926163725
$ws.Columns.Item(21).ColumnWidth = 6.755601926163725
